$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-18 (A:date serial, B:y0, C:y0_forecast, D:y1, E:y1_forecast)
$data = @(
    @(39765, 2008, $null, 2009, 2.741420124134053),
    @(40130, 2009, -3.872359107260159, 2010, -1.796159642965267),
    @(40494, 2010, 4.530477057343663, 2011, 3.372216600811506),
    @(40862, 2011, 6.833902841285977, 2012, 5.794553957309168),
    @(41228, 2012, 4.166536506645224, 2013, 4.453225076234824),
    @(41592, 2013, 2.669880057548091, 2014, 4.443665867214519),
    @(41957, 2014, 5.50293301232252, 2015, 4.93942156312468),
    @(42321, 2015, 4.829481320500406, 2016, 5.304093839766777),
    @(42689, 2016, 5.100281927437122, 2017, 4.773858406282372),
    @(43053, 2017, 5.161358932333737, 2018, 5.448823958443616),
    @(43418, 2018, 5.902681694119694, 2019, 5.03659417252571),
    @(43783, 2019, 3.884502719230132, 2020, 4.009670676786059),
    @(44159, 2020, -3.840397826549158, 2021, -0.236010050592228),
    @(44525, 2021, 0.4839811651348835, 2022, 2.177145583294293),
    @(44890, 2022, 2.06342951900429, 2023, 0.3896432785800652),
    @(45254, 2023, -2.156362896191677, 2024, -1.437335768580206),
    @(45618, 2024, -0.8205034771073372, 2025, 0.07297157746815053)
)

# Delete the last row (row 19), shrinking the table from 18 data rows to 17
$ws.Rows.Item(19).Delete()

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -eq $null) {
        $ws.Cells.Item($r, 3).Value = $null
    } else {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
